$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells to update: column -> row -> new text value
$updates = @{
    "D2" = "308.25"
    "E2" = "-0.13%"
    "D3" = "39.89"
    "E4" = "1.19%"
    "D5" = "0.08141"
    "E5" = "-0.28%"
    "D6" = "1.941"
    "E6" = "-2.45%"
    "D7" = "8.162"
    "E7" = "3.29%"
    "D8" = "4.240"
    "E8" = "1.38%"
    "D9" = "0.9295"
    "E9" = "-0.31%"
    "D10" = "0.1434"
    "E10" = "1.87%"
    "D11" = "0.1923"
    "E11" = "-1.31%"
    "D12" = "0.09087"
    "E12" = "-2.00%"
    "D13" = "0.03516"
    "E13" = "2.13%"
    "D14" = "0.09777"
    "E14" = "-0.78%"
    "D15" = "0.001395"
    "E15" = "-1.04%"
    "D16" = "0.005833"
    "E16" = "-4.24%"
    "D17" = "3.920"
    "E17" = "4.19%"
    "D18" = "3.380"
    "E18" = "-2.82%"
    "D19" = "0.3430"
    "E19" = "-0.56%"
    "D20" = "0.1313"
    "E20" = "0.73%"
    "D21" = "4.632"
    "E21" = "-3.81%"
    "E22" = "-1.19%"
    "D23" = "0.04374"
    "E23" = "-1.96%"
    "E24" = "-0.94%"
    "D25" = "0.004376"
    "E25" = "4.84%"
    "D26" = "0.0001300"
    "E26" = "-0.18%"
    "E27" = "-10.04%"
    "D39" = "0.02054"
    "E39" = "-3.17%"
    "D40" = "0.05071"
    "E40" = "-2.01%"
    "D41" = "0.007410"
    "E41" = "-0.65%"
    "D42" = "0.009850"
    "E42" = "-1.95%"
    "D43" = "0.1364"
    "E43" = "-0.45%"
    "D44" = "0.002130"
    "E44" = "-0.18%"
    "D45" = "0.009371"
    "E45" = "-3.25%"
    "D46" = "0.00006369"
    "E46" = "0.80%"
    "D47" = "0.00000000750"
    "E47" = "-0.08%"
    "D48" = "0.002714"
    "E49" = "-18.80%"
    "D50" = "0.00002101"
    "E50" = "-0.08%"
    "D51" = "0.0002001"
    "E51" = "-0.08%"
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    # Force text storage so numeric-looking / percent-looking strings
    # are not auto-converted to numbers by Excel.
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$addr]
    # Restore the default (unstyled) look so only the cell VALUE changes,
    # matching the source workbook which carries no explicit style on
    # these cells.
    $cell.Style = "Normal"
}
